$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.447.76'
$ws.Range("E2").Value = '  -1.68%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.649.75'
$ws.Range("E3").Value = '  -3.46%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.75'
$ws.Range("E5").Value = '  -0.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9993'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3652'
$ws.Range("E7").Value = '  -2.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '46.47'
$ws.Range("E8").Value = '  -6.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3244'
$ws.Range("E9").Value = '  -6.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.123'
$ws.Range("E10").Value = '  -7.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07014'
$ws.Range("E11").Value = '  -7.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9992'
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.963'
$ws.Range("E13").Value = '  -5.91%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.33'
$ws.Range("E14").Value = '  -9.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.603'
$ws.Range("E15").Value = '  -6.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.647.93'
$ws.Range("E16").Value = '  -3.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001039'
$ws.Range("E17").Value = '  -8.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06571'
$ws.Range("E18").Value = '  -2.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9985'
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '78.40'
$ws.Range("E20").Value = '  -7.65%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.917'
$ws.Range("E21").Value = '  -7.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.63'
$ws.Range("E22").Value = '  -10.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.54'
$ws.Range("E23").Value = '  -5.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.425.01'
$ws.Range("E24").Value = '  -1.83%  '
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.318'
$ws.Range("E26").Value = '  -17.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '146.62'
$ws.Range("E27").Value = '  -2.79%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.56'
$ws.Range("E28").Value = '  -9.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.830.89'
$ws.Range("E29").Value = '  -3.71%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.191'
$ws.Range("E30").Value = '  -4.81%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '123.95'
$ws.Range("E31").Value = '  -6.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.060'
$ws.Range("E32").Value = '  -4.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.693'
$ws.Range("E33").Value = '  -17.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08442'
$ws.Range("E34").Value = '  -4.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.657'
$ws.Range("E35").Value = '  -5.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.08'
$ws.Range("E36").Value = '  -13.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.201'
$ws.Range("E37").Value = '  -8.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.265'
$ws.Range("E38").Value = '  -1.40%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02232'
$ws.Range("E39").Value = '  -7.76%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06017'
$ws.Range("E40").Value = '  -9.91%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2062'
$ws.Range("E41").Value = '  -8.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.080'
$ws.Range("E42").Value = '  -13.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9988'
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5886'
$ws.Range("E44").Value = '  -9.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.763'
$ws.Range("E45").Value = '  -1.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.57'
$ws.Range("E46").Value = '  -9.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5607'
$ws.Range("E47").Value = '  -9.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.87'
$ws.Range("E48").Value = '  -6.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.942'
$ws.Range("E49").Value = '  -9.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06902'
$ws.Range("E50").Value = '  -5.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.179'
$ws.Range("E51").Value = '  -3.62%  '
